$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply per-cell updates derived from the canonical OOXML diff.
# D-column cells whose new value parses as a plain number need an explicit
# text format ("@") first, otherwise Excel COM auto-converts the assigned
# string into a numeric value (losing formatting / exact text, e.g. "1.00" -> 1).

$ws.Range('D2').Value = '68.785.99'
$ws.Range('E2').Value = '  +1.16%  '
$ws.Range('D3').Value = '3.280.29'
$ws.Range('E3').Value = '  +0.44%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '583.06'
$ws.Range('E5').Value = '  +0.29%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '185.78'
$ws.Range('E6').Value = '  +1.83%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -0.60%  '
$ws.Range('E9').Value = '  -0.40%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.65'
$ws.Range('E10').Value = '  -1.11%  '
$ws.Range('E11').Value = '  +1.16%  '
$ws.Range('D12').Value = '3.854.58'
$ws.Range('E12').Value = '  +0.44%  '
$ws.Range('E13').Value = '  -0.11%  '
$ws.Range('E14').Value = '  -0.38%  '
$ws.Range('D15').Value = '68.789.10'
$ws.Range('E15').Value = '  +1.09%  '
$ws.Range('E16').Value = '  +1.50%  '
$ws.Range('D17').Value = '3.268.18'
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '5.87'
$ws.Range('E18').Value = '  +0.28%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.63'
$ws.Range('E19').Value = '  +0.83%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '395.56'
$ws.Range('E20').Value = '  +5.10%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.73'
$ws.Range('E21').Value = '  +0.97%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '71.69'
$ws.Range('E22').Value = '  +0.60%  '
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('E24').Value = '  +1.34%  '
$ws.Range('E25').Value = '  +1.08%  '
$ws.Range('E26').Value = '  +4.24%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.74'
$ws.Range('E27').Value = '  +1.27%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.32%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('E30').Value = '  +0.33%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '23.12'
$ws.Range('E31').Value = '  +1.15%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.30'
$ws.Range('E32').Value = '  +2.51%  '
$ws.Range('B33').Value = 'Aptos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '7.17'
$ws.Range('E33').Value = '  +3.53%  '
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.51'
$ws.Range('E35').Value = '  -0.65%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '163.36'
$ws.Range('E36').Value = '  +0.71%  '
$ws.Range('E37').Value = '  +6.49%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.828'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '26.88'
$ws.Range('E39').Value = '  -0.16%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '4.62'
$ws.Range('E40').Value = '  -0.53%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '6.64'
$ws.Range('E41').Value = '  -2.87%  '
$ws.Range('E42').Value = '  -2.46%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '41.51'
$ws.Range('E43').Value = '  +1.55%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0691'
$ws.Range('E44').Value = '  +1.44%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '25.49'
$ws.Range('E45').Value = '  -0.70%  '
$ws.Range('D46').Value = '2.659.81'
$ws.Range('E46').Value = '  -1.02%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '342.91'
$ws.Range('E47').Value = '  -2.38%  '
$ws.Range('E48').Value = '  +0.81%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '6.37'
$ws.Range('E49').Value = '  +3.40%  '
$ws.Range('B50').Value = 'Arweave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '31.96'
$ws.Range('E50').Value = '  +2.25%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.994'
$ws.Range('E51').Value = '  -0.88%  '
